$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new column F ("estimate_se") before the old zval column ---
$ws.Columns("F:F").Insert()

# Header
$ws.Range("F1").Value = "estimate_se"

# Formula: concatenate estimate and se with a "±" separator
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 6).Formula = '=CONCATENATE(D' + $r + ',"±",E' + $r + ')'
}

# --- Sheet view tidy-up: drop the frozen/scrolled topLeftCell and select the new columns ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G2:H19").Select()
